# Tab2_Studienjahrgang.xlsx -- course-planning patch
# Expands the Studienjahrgang table from a single sample row to the full
# 21-row dataset (Bioanalytik/Zellbiologie, Chemie- und Bioprozesstechnik,
# Chemie, Medizininformatik, Medizintechnik, Pharmatechnologie, Umwelttechnik)
# and refreshes the AutoFilter / FilterDatabase range + selection to match.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet currently has header (row 1) + 1 data row (row 2) + 1 blank
# trailing row (row 3). Insert 20 rows above the blank trailing row so it
# becomes row 23, leaving rows 2-22 for the full dataset. Inserting copies
# row 2 formatting into the new rows and pushes the blank row down intact.
$ws.Rows("3:22").Insert()

$data = @(
    @("9322951", "2-L-B-LSBZ/18.a-SJ", "Bachelor of Science in Life Sciences Bioanalytik und Zellbiologie 2018", "9311263", "2-L-B-LSBZ/19"),
    @("9311732", "2-L-B-LSBZ/19.a-SJ", "Bachelor of Science in Life Sciences Bioanalytik und Zellbiologie 2019", "9311263", "2-L-B-LSBZ/19"),
    @("9369191", "2-L-B-LSBZ/20.a-SJ", "Bachelor of Science in Life Sciences Bioanalytik und Zellbiologie 2020", "9311263", "2-L-B-LSBZ/19"),
    @("9322952", "2-L-B-LSCB/18.a-SJ", "Bachelor of Science in Life Sciences Chemie- und Bioprozesstechnik 2018", "9311288", "2-L-B-LSCB/19"),
    @("9311734", "2-L-B-LSCB/19.a-SJ", "Bachelor of Science in Life Sciences Chemie- und Bioprozesstechnik 2019", "9311288", "2-L-B-LSCB/19"),
    @("9369194", "2-L-B-LSCB/20.a-SJ", "Bachelor of Science in Life Sciences Chemie- und Bioprozesstechnik 2020", "9311288", "2-L-B-LSCB/19"),
    @("9322953", "2-L-B-LSCH/18.a-SJ", "Bachelor of Science in Life Sciences Chemie 2018", "9311212", "2-L-B-LSCH/19"),
    @("9311735", "2-L-B-LSCH/19.a-SJ", "Bachelor of Science in Life Sciences Chemie 2019", "9311212", "2-L-B-LSCH/19"),
    @("9369195", "2-L-B-LSCH/20.a-SJ", "Bachelor of Science in Life Sciences Chemie 2020", "9311212", "2-L-B-LSCH/19"),
    @("9322955", "2-L-B-LSMI/18.a-SJ", "Bachelor of Science in Life Sciences Medizininformatik 2018", "9304956", "2-L-B-LSMI/19"),
    @("9311737", "2-L-B-LSMI/19.a-SJ", "Bachelor of Science in Life Sciences Medizininformatik 2019", "9304956", "2-L-B-LSMI/19"),
    @("9369197", "2-L-B-LSMI/20.a-SJ", "Bachelor of Science in Life Sciences Medizininformatik 2020", "9304956", "2-L-B-LSMI/19"),
    @("9322956", "2-L-B-LSMT/18.a-SJ", "Bachelor of Science in Life Sciences Medizintechnik 2018", "9309072", "2-L-B-LSMT/19"),
    @("9311738", "2-L-B-LSMT/19.a-SJ", "Bachelor of Science in Life Sciences Medizintechnik 2019", "9309072", "2-L-B-LSMT/19"),
    @("9369169", "2-L-B-LSMT/20.a-SJ", "Bachelor of Science in Life Sciences Medizintechnik 2020", "9309072", "2-L-B-LSMT/19"),
    @("9322957", "2-L-B-LSPT/18.a-SJ", "Bachelor of Science in Life Sciences Pharmatechnologie 2018", "9311111", "2-L-B-LSPT/19"),
    @("9311739", "2-L-B-LSPT/19.a-SJ", "Bachelor of Science in Life Sciences Pharmatechnologie 2019", "9311111", "2-L-B-LSPT/19"),
    @("9369199", "2-L-B-LSPT/20.a-SJ", "Bachelor of Science in Life Sciences Pharmatechnologie 2020", "9311111", "2-L-B-LSPT/19"),
    @("9322950", "2-L-B-LSUT/18.a-SJ", "Bachelor of Science in Life Sciences Umwelttechnik 2018", "9310716", "2-L-B-LSUT/19"),
    @("9311740", "2-L-B-LSUT/19.a-SJ", "Bachelor of Science in Life Sciences Umwelttechnik 2019", "9310716", "2-L-B-LSUT/19"),
    @("9369200", "2-L-B-LSUT/20.a-SJ", "Bachelor of Science in Life Sciences Umwelttechnik 2020", "9310716", "2-L-B-LSUT/19"),
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = 2 + $i
    $rec = $data[$i]
    $ws.Cells.Item($row, 1).Value = [double]$rec[0]
    $ws.Cells.Item($row, 2).Value = $rec[1]
    $ws.Cells.Item($row, 3).Value = $rec[2]
    $ws.Cells.Item($row, 4).Value = [double]$rec[3]
    $ws.Cells.Item($row, 5).Value = $rec[4]
}

# Re-apply the AutoFilter over the grown range A1:E22 (toggle off first --
# calling .AutoFilter() while a filter is already active just clears it).
$ws.AutoFilterMode = $false
[void]$ws.Range("A1:E22").AutoFilter()

# Keep the workbook-level hidden _FilterDatabase name in sync with the
# new filter range, same as Excel does automatically when a live user
# resizes an AutoFilter.
$names = $wb.Names
for ($i = 1; $i -le $names.Count; $i++) {
    $n = $names.Item($i)
    if ($n.Name -eq "Sheet!_FilterDatabase") {
        $n.RefersTo = "=Sheet!`$A`$1:`$E`$22"
    }
}

# Match the final selection recorded in the sheet: column B across the
# data rows.
[void]$ws.Range("B2:B22").Select()

